$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their text formatting (avoid Excel auto-converting
# numeric-looking strings like "1.00" into the number 1). We use the classic
# leading-apostrophe 'treat as text' prefix that Excel understands natively.

$ws.Range("D2").Value = "'59.292.72"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "'3.180.29"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'533.02"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'142.33"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  +11.54%  "
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'0.441"
$ws.Range("E10").Value = "  +6.84%  "
$ws.Range("D11").Value = "'0.113"
$ws.Range("E11").Value = "  +4.73%  "
$ws.Range("D12").Value = "'3.729.78"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "'25.96"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("E15").Value = "  +4.47%  "
$ws.Range("D16").Value = "'59.310.79"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D17").Value = "'3.208.31"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "'6.25"
$ws.Range("E18").Value = "  +2.57%  "
$ws.Range("D19").Value = "'13.04"
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").Value = "'8.19"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "'377.03"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("D24").Value = "'69.83"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'8.43"
$ws.Range("E27").Value = "  +15.25%  "
$ws.Range("D28").Value = "'0.0₃0874"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").Value = "'22.46"
$ws.Range("E29").Value = "  +4.79%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'6.06"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").Value = "'5.24"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  +4.40%  "
$ws.Range("D35").Value = "'157.33"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("D37").Value = "'0.0718"
$ws.Range("E37").Value = "  +6.82%  "
$ws.Range("D38").Value = "'25.44"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'2.711.43"
$ws.Range("E39").Value = "  +7.70%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0294"
$ws.Range("E42").Value = "  +8.91%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.726"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").Value = "'39.18"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'3.221.34"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  +12.54%  "
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").Value = "'20.34"
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("D51").Value = "'0.761"
$ws.Range("E51").Value = "  +1.91%  "
